# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    3  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    4  = @{ B = 3.230985683306322;  C = 10.29869402782916;  D = 3.900430680208489;  E = 8.660232485948974; G = 26.09034287729295 }
    5  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 3.781711156805759 }
    6  = @{ B = 3.230985683306322;  C = 0.3127903958511391; D = 3.900430680208489;  E = 8.660232485948974; G = 16.10443924531492 }
    7  = @{ B = 0.003994804209775715; C = 0.3127903958511391; D = 0.8054896365839992; E = 0.496779210170732; G = 1.619054046815646 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    9  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 4.429675500412797 }
    10 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    11 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 6.740334628841572 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
}
